$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 24999
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 24999
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 24999
$ws.Range("N7").Value = -25223

$ws.Range("H14").Value = 24999
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 24999
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 24999
$ws.Range("N14").Value = -25381

$ws.Range("H17").Value = 2024.75
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 2366.3333
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 7098.999899999999
$ws.Range("M17").Value = -2832
$ws.Range("N17").Value = -7434.999899999999

$ws.Range("H19").Value = 530.619
$ws.Range("I19").Value = 644.1
$ws.Range("J19").Value = 427.45456
$ws.Range("K19").Value = 644.1
$ws.Range("L19").Value = 427.45456
$ws.Range("M19").Value = -469.1
$ws.Range("N19").Value = -777.45456

$ws.Range("H41").Value = 706.9259
$ws.Range("I41").Value = 596.44446
$ws.Range("J41").Value = 927.8889
$ws.Range("K41").Value = 596.44446
$ws.Range("L41").Value = 927.8889
$ws.Range("M41").Value = -156.44446
$ws.Range("N41").Value = -1807.8889

$ws.Range("H43").Value = 2723
$ws.Range("I43").Value = 2100
$ws.Range("J43").Value = 3969
$ws.Range("K43").Value = 2100
$ws.Range("L43").Value = 3969
$ws.Range("M43").Value = -2031

$ws.Range("H53").Value = 353.69232
$ws.Range("I53").Value = 145.85715
$ws.Range("J53").Value = 596.1667
$ws.Range("K53").Value = 145.85715
$ws.Range("L53").Value = 596.1667
$ws.Range("M53").Value = 491.14285
$ws.Range("N53").Value = -1870.1667

$ws.Range("H70").Value = 4504.467
$ws.Range("I70").Value = 3667
$ws.Range("J70").Value = 4713.8335
$ws.Range("K70").Value = 11001
$ws.Range("L70").Value = 14141.5005
$ws.Range("M70").Value = -10731
$ws.Range("N70").Value = -14681.5005

$ws.Range("H73").Value = 4504.467
$ws.Range("I73").Value = 3667
$ws.Range("J73").Value = 4713.8335
$ws.Range("K73").Value = 11001
$ws.Range("L73").Value = 14141.5005
$ws.Range("M73").Value = -10065
$ws.Range("N73").Value = -16013.5005

$ws.Range("H132").Value = 18388.621
$ws.Range("I132").Value = 1052.7307
$ws.Range("J132").Value = 168633
$ws.Range("K132").Value = 3158.1921
$ws.Range("L132").Value = 505899
$ws.Range("M132").Value = -628.1921000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 30000
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 30000
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 30000
$ws.Range("N34").Value = -30542

$ws.Range("H61").Value = 5135.364
$ws.Range("I61").Value = 4346.0435
$ws.Range("J61").Value = 6950.8
$ws.Range("K61").Value = 4346.0435
$ws.Range("L61").Value = 6950.8
$ws.Range("M61").Value = -4134.0435

$ws.Range("H128").Value = 60214.5
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 60214.5
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 60214.5
$ws.Range("N128").Value = -70174.5

$ws.Range("H132").Value = 5241.974
$ws.Range("I132").Value = 5088.533
$ws.Range("J132").Value = 5753.4443
$ws.Range("K132").Value = 15265.599
$ws.Range("L132").Value = 17260.3329
$ws.Range("M132").Value = -12735.599
$ws.Range("N132").Value = -22320.3329

$ws.Range("H136").Value = 5135.364
$ws.Range("I136").Value = 4346.0435
$ws.Range("J136").Value = 6950.8
$ws.Range("K136").Value = 13038.1305
$ws.Range("L136").Value = 20852.4
$ws.Range("M136").Value = -10488.1305

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4792.6665
$ws.Range("I105").Value = 2841.8572
$ws.Range("J105").Value = 6499.625
$ws.Range("K105").Value = 2841.8572
$ws.Range("L105").Value = 6499.625
$ws.Range("M105").Value = -1094.8572
$ws.Range("N105").Value = -9993.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 392
$ws.Range("I7").Value = 180
$ws.Range("J7").Value = 445
$ws.Range("K7").Value = 180
$ws.Range("L7").Value = 445
$ws.Range("M7").Value = -67
$ws.Range("N7").Value = -671

$ws.Range("H31").Value = 2937.0815
$ws.Range("I31").Value = 2278.8572
$ws.Range("J31").Value = 4582.643
$ws.Range("K31").Value = 2278.8572
$ws.Range("L31").Value = 4582.643
$ws.Range("M31").Value = -1983.8572

$ws.Range("H34").Value = 2937.0815
$ws.Range("I34").Value = 2278.8572
$ws.Range("J34").Value = 4582.643
$ws.Range("K34").Value = 2278.8572
$ws.Range("L34").Value = 4582.643
$ws.Range("M34").Value = -2076.8572

$ws.Range("H125").Value = 67500
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 67500
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 67500
$ws.Range("N125").Value = -72420

$ws.Range("H134").Value = 3948.4814
$ws.Range("I134").Value = 2621.6843
$ws.Range("J134").Value = 7099.625
$ws.Range("K134").Value = 7865.0529
$ws.Range("L134").Value = 21298.875
$ws.Range("M134").Value = -5330.0529

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 324.2963
$ws.Range("I107").Value = 392
$ws.Range("J107").Value = 315.83334
$ws.Range("K107").Value = 1176
$ws.Range("L107").Value = 947.5000200000001
$ws.Range("M107").Value = 744
$ws.Range("N107").Value = -4787.50002

$ws.Range("H122").Value = 1969.2941
$ws.Range("I122").Value = 1779.8
$ws.Range("J122").Value = 2048.25
$ws.Range("K122").Value = 16018.2
$ws.Range("L122").Value = 18434.25
$ws.Range("M122").Value = -13568.2
$ws.Range("N122").Value = -23334.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 470
$ws.Range("I22").Value = 470
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 470
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 59
$ws.Range("N22").ClearContents()

$ws.Range("H133").Value = 70000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 70000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5400
$ws.Range("I7").Value = 3406.6
$ws.Range("J7").Value = 6396.7
$ws.Range("K7").Value = 3406.6
$ws.Range("L7").Value = 6396.7
$ws.Range("M7").Value = -3294.6

$ws.Range("H55").Value = 2113.25
$ws.Range("I55").Value = 3564
$ws.Range("J55").Value = 1242.8
$ws.Range("K55").Value = 3564
$ws.Range("L55").Value = 1242.8
$ws.Range("M55").Value = -3391

$ws.Range("H126").Value = 5400
$ws.Range("I126").Value = 3406.6
$ws.Range("J126").Value = 6396.7
$ws.Range("K126").Value = 10219.8
$ws.Range("L126").Value = 19190.1
$ws.Range("M126").Value = -7749.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 15000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 15000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 15000
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -15826

$ws.Range("H81").Value = 9806458
$ws.Range("I81").Value = 2073.72
$ws.Range("J81").Value = 37040860
$ws.Range("K81").Value = 4147.44
$ws.Range("L81").Value = 74081720
$ws.Range("M81").Value = -3086.44

$ws.Range("H84").Value = 9806458
$ws.Range("I84").Value = 2073.72
$ws.Range("J84").Value = 37040860
$ws.Range("K84").Value = 20737.2
$ws.Range("L84").Value = 370408600
$ws.Range("M84").Value = -15433.2

$ws.Range("H113").Value = 632.35297
$ws.Range("I113").Value = 410.16
$ws.Range("J113").Value = 1249.5555
$ws.Range("K113").Value = 1230.48
$ws.Range("L113").Value = 3748.6665
$ws.Range("M113").Value = 939.52
$ws.Range("N113").Value = -8088.666499999999
